$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Drawdown_Peak" test row (row 85) to "Drawdowns"
$ws.Range("A85").Value = "Drawdowns1"
$ws.Range("B85").Value = "Test drawdowns"
$ws.Range("C85").Value = "Drawdowns_test1"

# Update the active selection to match the committed state
$ws.Range("F84").Select()
